# The commit swaps the presentation's theme ("Integral") for the
# default "Office Theme" palette (the 25 May 2020 -> 08 Jun 2020 edit
# re-themed the deck). Concretely the 12 theme colours (dk1/lt1/dk2/lt2/
# accent1-6/hlink/folHlink) defined on the slide master's colour scheme
# change from the "Integral" values to the stock "Office Theme" values.
#
# PowerPoint's ColorScheme.Colors(n).RGB is BGR-packed (0xBBGGRR), so the
# target RGB hex values below are byte-swapped before assignment.

$p  = $ppt.ActivePresentation
$sm = $p.SlideMaster
$cs = $sm.ColorScheme

$cs.Colors(1).RGB  = 0x000000  # dk1      -> 000000
$cs.Colors(2).RGB  = 0xFFFFFF  # lt1      -> FFFFFF
$cs.Colors(3).RGB  = 0x6A5444  # dk2      -> 44546A
$cs.Colors(4).RGB  = 0xE6E6E7  # lt2      -> E7E6E6
$cs.Colors(5).RGB  = 0xD59B5B  # accent1  -> 5B9BD5
$cs.Colors(6).RGB  = 0x317DED  # accent2  -> ED7D31
$cs.Colors(7).RGB  = 0xA5A5A5  # accent3  -> A5A5A5
$cs.Colors(8).RGB  = 0x00C0FF  # accent4  -> FFC000
$cs.Colors(9).RGB  = 0xC47244  # accent5  -> 4472C4
$cs.Colors(10).RGB = 0x47AD70  # accent6  -> 70AD47
$cs.Colors(11).RGB = 0xC16305  # hlink    -> 0563C1
$cs.Colors(12).RGB = 0x724F95  # folHlink -> 954F72
